{"js": "// Remove the \"\u00a9 2020 ... Creative Commons Attribution\" copyright paragraph\n// together with the two empty paragraphs that immediately precede it\n// (a plain blank paragraph and a blank page-break paragraph). The\n// paragraph right after the copyright notice (another blank paragraph)\n// is left untouched, so it now follows directly after the\n// \"LOQ4095: ...\" paragraph.\n\nconst body = context.document.body;\n\n// Locate the copyright paragraph by its distinctive leading text.\nconst results = body.search(\"\u00a9 2020\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the \"\u00a9 2020\" copyright paragraph.');\n}\n\nconst copyrightPara = results.items[0].paragraphs.getFirst();\nconst pageBreakBlankPara = copyrightPara.getPrevious(); // blank paragraph with pageBreakBefore\nconst plainBlankPara = pageBreakBlankPara.getPrevious(); // plain blank paragraph\n\n// Delete the three paragraphs (order doesn't matter for correctness,\n// but deleting the anchors we already captured is safe either way).\ncopyrightPara.delete();\npageBreakBlankPara.delete();\nplainBlankPara.delete();\n\nawait context.sync();\n", "ps1": "# Remove the \"\u00a9 2020 ... Creative Commons Attribution\" copyright paragraph\n# together with the two empty paragraphs that immediately precede it (a\n# plain blank paragraph and a blank page-break paragraph). The paragraph\n# right after the copyright notice (another blank paragraph) is left\n# untouched, so it now follows directly after the \"LOQ4095: ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the copyright paragraph. Search on an ASCII-safe substring of its\n# text (avoids any \"\u00a9\" transcoding quirks) rather than a hard-coded index.\n$count = $d.Paragraphs.Count\n$foundIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Contact: luizeleno*\") {\n        $foundIndex = $i\n        break\n    }\n}\n\nif ($foundIndex -eq -1) {\n    throw \"Could not find the copyright paragraph.\"\n}\n\n# Delete highest index first so the lower indices stay valid.\n$d.Paragraphs.Item($foundIndex).Range.Delete()\n$d.Paragraphs.Item($foundIndex - 1).Range.Delete()\n$d.Paragraphs.Item($foundIndex - 2).Range.Delete()\n"}
